$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.05
$ws.Range("I3").Value = 4.5
$ws.Range("K3").Value = 1.83
$ws.Range("W3").Value = 7
$ws.Range("X3").Value = 1.1
$ws.Range("AS3").Value = 67
$ws.Range("G4").Value = 2.55
$ws.Range("H4").Value = 2.9
$ws.Range("I4").Value = 3.2
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("AE4").Value = 11
$ws.Range("AF4").Value = 26
$ws.Range("AJ4").Value = 5.5
$ws.Range("I5").Value = 2.88
$ws.Range("O5").Value = 1.8
$ws.Range("P5").Value = 1.91
$ws.Range("S5").Value = 3.6
$ws.Range("T5").Value = 1.29
$ws.Range("U5").Value = 6.8
$ws.Range("AD5").Value = 12
$ws.Range("AI5").Value = 4.75
$ws.Range("M6").Value = 1.14
$ws.Range("N6").Value = 5.5
$ws.Range("Z6").Value = 2
$ws.Range("S7").Value = 3.1
$ws.Range("T7").Value = 1.36
$ws.Range("W7").Value = 6.5
$ws.Range("X7").Value = 1.11
$ws.Range("N8").Value = 9
$ws.Range("O8").Value = 1.33
$ws.Range("P8").Value = 3.25
$ws.Range("S8").Value = 2.1
$ws.Range("T8").Value = 1.7
$ws.Range("G9").Value = 2.3
$ws.Range("I9").Value = 3.25
$ws.Range("J9").Value = 3.1
$ws.Range("L9").Value = 3.75
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8
$ws.Range("AA9").Value = 1.83
$ws.Range("AB9").Value = 1.83
$ws.Range("AG9").Value = 21
$ws.Range("AO9").Value = 15
$ws.Range("G10").Value = 2.1
$ws.Range("H10").Value = 3.25
$ws.Range("I10").Value = 3.2
$ws.Range("J10").Value = 2.7
$ws.Range("K10").Value = 2.15
$ws.Range("L10").Value = 3.7
$ws.Range("Y10").Value = 1.39
$ws.Range("Z10").Value = 2.77
$ws.Range("AF10").Value = 19.5
$ws.Range("AG10").Value = 17.5
$ws.Range("AQ10").Value = 45
$ws.Range("AR10").Value = 29
$ws.Range("M11").Value = 1.05
$ws.Range("O11").Value = 1.25
$ws.Range("S11").Value = 1.88
$ws.Range("T11").Value = 1.98
$ws.Range("X11").Value = 1.33
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("O12").Value = 1.53
$ws.Range("P12").Value = 2.38
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 1.85
$ws.Range("X12").Value = 1.14
$ws.Range("G13").Value = 2.55
$ws.Range("I13").Value = 2.63
$ws.Range("J13").Value = 3.25
$ws.Range("M13").Value = 1.06
$ws.Range("O13").Value = 1.33
$ws.Range("X13").Value = 1.25
$ws.Range("AC13").Value = 8
$ws.Range("AD13").Value = 12
$ws.Range("AG13").Value = 21
$ws.Range("AP13").Value = 11
